# "Generate Report for Handback" — the handback files have now landed in
# sync with en-US, so populate the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns for the 125b7fa1-... entry and
# flip its status, on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$sheetInfo = @{
    "zh-cn" = @{ HandbackTime = "2016-03-10 16:34:26" };
    "de-de" = @{ HandbackTime = "2016-03-10 16:34:35" };
}

foreach ($sheetName in $sheetInfo.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetInfo[$sheetName]

    # Values already present on row 2 (the 125b7fa1-ae8f-... file) — re-use
    # the handoff file name / source file name for the handback columns,
    # since the handback content matches what was handed off.
    $handoffDisplay = $ws.Range("A2").Text
    $handoffFileDisplay = $ws.Range("C2").Text

    # Row 2 (125b7fa1-ae8f-46ad-ab80-0ed7a0f66097): now in sync with en-US —
    # mirror the handoff file/date into the "Latest Target File" /
    # "Latest Handback File" columns and record when the handback landed.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("E2").Value = $handoffDisplay
    $ws.Range("F2").Value = $handoffFileDisplay
    $ws.Range("G2").Value = $info.HandbackTime
}

# Mirror the same hyperlink that's on A2/C2 onto the newly-populated E2/F2
# cells (Latest Target File / Latest Handback File) for each language sheet.
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11fc11f582a6164b74e3c018efbd82d5247f3f3e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf") | Out-Null

$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/468abc9da3edc8659891fe1397da2d14b2e02887/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf") | Out-Null
